$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 21
$ws.Range("A9").Value = 32
$ws.Range("A10").Value = 13

$ws.Range("A8").Select()
